# fix(docx): fix OOXMLValidator error on KeywordTok output
#
# The <w:rPr> children of several custom character styles in styles.xml
# were serialized with <w:color> before <w:b>/<w:i>, which does not match
# the sequence required by wml.xsd (CT_RPr expects b/bCs/i/iCs/... before
# color). Re-apply bold/italic on each affected style so the document
# model re-emits the run properties in the schema-mandated order
# (b, i, color).

$d = $word.ActiveDocument

# Styles whose <w:rPr> only contains <w:b/> and <w:color/> (wrong order:
# color before b).
$boldOnly = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

# Styles whose <w:rPr> only contains <w:i/> and <w:color/> (wrong order:
# color before i).
$italicOnly = @(
    "CommentTok",
    "DocumentationTok"
)

# Styles whose <w:rPr> contains <w:b/>, <w:i/> and <w:color/> (wrong
# order: color before b/i).
$boldItalic = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($name in $boldOnly) {
    $s = $d.Styles($name)
    $s.Font.Bold = 0
    $s.Font.Bold = -1
}

foreach ($name in $italicOnly) {
    $s = $d.Styles($name)
    $s.Font.Italic = 0
    $s.Font.Italic = -1
}

foreach ($name in $boldItalic) {
    $s = $d.Styles($name)
    $s.Font.Bold = 0
    $s.Font.Bold = -1
    $s.Font.Italic = 0
    $s.Font.Italic = -1
}
